$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 255; existing rows 255-304 shift down to 256-305,
# inheriting formatting (incl. the date style on column D) from the row
# that used to occupy 255.
$ws.Rows.Item(255).Insert()

# Populate the newly inserted row 255 with the new weekly record.
$ws.Cells.Item(255, 1).Value2 = 6
$ws.Cells.Item(255, 2).Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(255, 3).Value2 = "Metropolitana"
$ws.Cells.Item(255, 4).Value2 = 44995
$ws.Cells.Item(255, 5).Value2 = 13
$ws.Cells.Item(255, 6).Value2 = 100112029
$ws.Cells.Item(255, 7).Value2 = "Orégano"
$ws.Cells.Item(255, 8).Value2 = "Sin especificar"
$ws.Cells.Item(255, 9).Value2 = "Primera"
$ws.Cells.Item(255, 10).Value2 = 42
$ws.Cells.Item(255, 11).Value2 = 19000
$ws.Cells.Item(255, 12).Value2 = 20000
$ws.Cells.Item(255, 13).Value2 = 19476
$ws.Cells.Item(255, 14).Value2 = "$/docena de atados"
$ws.Cells.Item(255, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(255, 16).Value2 = 6492
$ws.Cells.Item(255, 17).Value2 = 3
$ws.Cells.Item(255, 18).Value2 = "Hortaliza"
